$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.414949178695679
$ws.Range("B1").Value = 1.968210458755493
$ws.Range("C1").Value = 2.408133745193481
$ws.Range("D1").Value = 4.812595844268799
$ws.Range("E1").Value = 0.916772723197937
